$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.981.18"
$ws.Range("E2").Value = "  -3.74%  "
$ws.Range("D3").Value = "3.132.79"
$ws.Range("E3").Value = "  -5.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -6.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.40"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D8").Value = "3.130.69"
$ws.Range("E8").Value = "  -5.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.445"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.24"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.110"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -7.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -5.07%  "
$ws.Range("D13").Value = "3.659.91"
$ws.Range("E13").Value = "  -5.51%  "
$ws.Range("E14").Value = "  -2.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.66%  "
$ws.Range("D16").Value = "3.126.40"
$ws.Range("E16").Value = "  -5.18%  "
$ws.Range("D17").Value = "57.825.56"
$ws.Range("E17").Value = "  -4.01%  "
$ws.Range("E18").Value = "  -7.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.04"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.02"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -6.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.78"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.85%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.79"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -7.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.505"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -6.71%  "
$ws.Range("D26").Value = "3.243.75"
$ws.Range("E26").Value = "  -5.86%  "
$ws.Range("B27").Value = "Kaspa"
$ws.Range("C27").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0957"
$ws.Range("E28").Value = "  -6.83%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.81"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.87"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -8.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.91"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.29%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.24"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "21.58"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.32"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.80"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -6.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.20"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.38"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.49"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0692"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -6.13%  "
$ws.Range("D42").Value = "3.160.70"
$ws.Range("E42").Value = "  -5.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.23"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.684"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -9.08%  "
$ws.Range("B45").Value = "ONDO"
$ws.Range("C45").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.07"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.00%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.92"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.43%  "
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.46"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -8.19%  "
$ws.Range("D49").Value = "2.253.06"
$ws.Range("E49").Value = "  -4.94%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.19"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.64%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.11%  "
